$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 83334856
$ws.Range("J29").Value = 4042
$ws.Range("L29").Value = 12126
$ws.Range("N29").Value = -12688
$ws.Range("H38").Value = 9337.15
$ws.Range("I38").Value = 13205.556
$ws.Range("K38").Value = 39616.66800000001
$ws.Range("M38").Value = -39244.66800000001
$ws.Range("H100").Value = 5434
$ws.Range("J100").Value = 7527.1113
$ws.Range("L100").Value = 7527.1113
$ws.Range("N100").Value = -8609.1113
$ws.Range("H116").Value = 5251.25
$ws.Range("I116").Value = 5251.25
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 5251.25
$ws.Range("L116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("N116").Value = -1809.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 129.16667
$ws.Range("I5").Value = 54.666668
$ws.Range("K5").Value = 54.666668
$ws.Range("M5").Value = 57.333332
$ws.Range("H32").Value = 3472.5083
$ws.Range("I32").Value = 2912.8965
$ws.Range("J32").Value = 14291.667
$ws.Range("K32").Value = 2912.8965
$ws.Range("L32").Value = 14291.667
$ws.Range("M32").Value = -2625.8965
$ws.Range("N32").Value = -14865.667
$ws.Range("H63").Value = 10495.889
$ws.Range("I63").Value = 12997.667
$ws.Range("J63").Value = 9995.532999999999
$ws.Range("K63").Value = 12997.667
$ws.Range("L63").Value = 9995.532999999999
$ws.Range("M63").Value = -12311.667
$ws.Range("N63").Value = -11367.533
$ws.Range("H66").Value = 10495.889
$ws.Range("I66").Value = 12997.667
$ws.Range("J66").Value = 9995.532999999999
$ws.Range("K66").Value = 64988.335
$ws.Range("L66").Value = 49977.66499999999
$ws.Range("M66").Value = -61556.335
$ws.Range("N66").Value = -56841.66499999999
$ws.Range("H74").Value = 1920.0952
$ws.Range("I74").Value = 1906.1
$ws.Range("K74").Value = 1906.1
$ws.Range("M74").Value = -1032.1
$ws.Range("H77").Value = 1920.0952
$ws.Range("I77").Value = 1906.1
$ws.Range("K77").Value = 9530.5
$ws.Range("M77").Value = -5162.5
$ws.Range("H102").Value = 3415
$ws.Range("I102").Value = 2454.1
$ws.Range("K102").Value = 2454.1
$ws.Range("M102").Value = -832.0999999999999
$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("L110").ClearContents()
$ws.Range("M110").ClearContents()
$ws.Range("N110").Value = 0
$ws.Range("H111").Value = 43000
$ws.Range("I111").Value = 43000
$ws.Range("K111").Value = 43000
$ws.Range("M111").Value = -38910
$ws.Range("H124").Value = 44192.57
$ws.Range("J124").Value = 44391.332
$ws.Range("L124").Value = 44391.332
$ws.Range("N124").Value = -54211.332
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").ClearContents()
$ws.Range("N131").Value = 0
$ws.Range("H132").Value = 3100
$ws.Range("I132").Value = 3100
$ws.Range("K132").Value = 9300
$ws.Range("M132").Value = -6770

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 129.16667
$ws.Range("I4").Value = 54.666668
$ws.Range("K4").Value = 54.666668
$ws.Range("M4").Value = 60.333332
$ws.Range("H107").Value = 6605.0835
$ws.Range("I107").Value = 7087
$ws.Range("K107").Value = 7087
$ws.Range("M107").Value = -5167
$ws.Range("H126").Value = 57500
$ws.Range("J126").Value = 57500
$ws.Range("L126").Value = 57500
$ws.Range("N126").Value = -67380
$ws.Range("H127").Value = 115000
$ws.Range("J127").Value = 115000
$ws.Range("L127").Value = 115000
$ws.Range("N127").Value = -124920
$ws.Range("H132").Value = 74493.836
$ws.Range("J132").Value = 74493.836
$ws.Range("L132").Value = 74493.836
$ws.Range("N132").Value = -84613.836

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2750.2856
$ws.Range("I99").Value = 3150.4
$ws.Range("J99").Value = 1750
$ws.Range("K99").Value = 3150.4
$ws.Range("L99").Value = 1750
$ws.Range("M99").Value = -1652.4
$ws.Range("N99").Value = -4746
$ws.Range("H122").Value = 2000
$ws.Range("I122").Value = 2000
$ws.Range("K122").Value = 6000
$ws.Range("M122").Value = -3550
$ws.Range("H126").Value = 2750.2856
$ws.Range("I126").Value = 3150.4
$ws.Range("J126").Value = 1750
$ws.Range("K126").Value = 9451.200000000001
$ws.Range("L126").Value = 5250
$ws.Range("M126").Value = -6981.200000000001
$ws.Range("N126").Value = -10190
$ws.Range("H132").Value = 2264.0908
$ws.Range("I132").Value = 2100.6667
$ws.Range("K132").Value = 6302.000100000001
$ws.Range("M132").Value = -3772.000100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("H122").Value = 3224.7144
$ws.Range("J122").Value = 3224.7144
$ws.Range("L122").Value = 29022.4296
$ws.Range("N122").Value = -33922.4296

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 19499.5
$ws.Range("I102").Value = 9000
$ws.Range("J102").Value = 29999
$ws.Range("K102").Value = 9000
$ws.Range("L102").Value = 29999
$ws.Range("M102").Value = -7378
$ws.Range("N102").Value = -33243
$ws.Range("H104").Value = 41987.5
$ws.Range("J104").Value = 41987.5
$ws.Range("L104").Value = 41987.5
$ws.Range("N104").Value = -48975.5
$ws.Range("H122").Value = 44873.75
$ws.Range("I122").Value = 49831.832
$ws.Range("K122").Value = 149495.496
$ws.Range("M122").Value = -147045.496
$ws.Range("H126").Value = 2947
$ws.Range("I126").Value = 3100
$ws.Range("J126").Value = 2794
$ws.Range("K126").Value = 9300
$ws.Range("L126").Value = 8382
$ws.Range("M126").Value = -6830
$ws.Range("N126").Value = -13322
$ws.Range("H132").Value = 3050
$ws.Range("I132").Value = 3050
$ws.Range("K132").Value = 9150
$ws.Range("M132").Value = -6620
$ws.Range("H141").Value = 80000
$ws.Range("J141").Value = 80000
$ws.Range("L141").Value = 80000
$ws.Range("N141").Value = -90360

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1433
$ws.Range("I16").Value = 470.7143
$ws.Range("J16").Value = 2395.2856
$ws.Range("K16").Value = 470.7143
$ws.Range("L16").Value = 2395.2856
$ws.Range("M16").Value = -300.7143
$ws.Range("N16").Value = -2735.2856
$ws.Range("H40").Value = 4455.75
$ws.Range("I40").Value = 3008.4736
$ws.Range("J40").Value = 7511.1113
$ws.Range("K40").Value = 3008.4736
$ws.Range("L40").Value = 7511.1113
$ws.Range("M40").Value = -2872.4736
$ws.Range("N40").Value = -7783.1113
$ws.Range("H46").Value = 1727.825
$ws.Range("I46").Value = 1254.6428
$ws.Range("J46").Value = 1982.6154
$ws.Range("K46").Value = 1254.6428
$ws.Range("L46").Value = 1982.6154
$ws.Range("M46").Value = -1066.6428
$ws.Range("N46").Value = -2358.6154
$ws.Range("H55").Value = 1258.7222
$ws.Range("I55").Value = 204.76923
$ws.Range("K55").Value = 204.76923
$ws.Range("M55").Value = -31.76922999999999
$ws.Range("H93").Value = 6575
$ws.Range("J93").Value = 7777.778
$ws.Range("L93").Value = 7777.778
$ws.Range("N93").Value = -10273.778
$ws.Range("H103").Value = 24249.5
$ws.Range("J103").Value = 24249.5
$ws.Range("L103").Value = 24249.5
$ws.Range("N103").Value = -26593.5
$ws.Range("H132").Value = 3648.3215
$ws.Range("I132").Value = 3437.6667
$ws.Range("K132").Value = 10313.0001
$ws.Range("M132").Value = -7783.000100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1629
$ws.Range("I113").Value = 1253.1111
$ws.Range("J113").Value = 2051.875
$ws.Range("K113").Value = 3759.3333
$ws.Range("L113").Value = 6155.625
$ws.Range("M113").Value = -1589.3333
$ws.Range("N113").Value = -10495.625
